$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1730.25
$ws.Range("J29").Value = 2253.3333
$ws.Range("L29").Value = 6759.999899999999
$ws.Range("N29").Value = -7321.999899999999
$ws.Range("H32").Value = 1397.8
$ws.Range("H38").Value = 1626.3
$ws.Range("J38").Value = 1978.4062
$ws.Range("L38").Value = 5935.2186
$ws.Range("N38").Value = -6679.2186
$ws.Range("H39").Value = 76.8
$ws.Range("I39").Value = 72.25
$ws.Range("J39").Value = 95
$ws.Range("K39").Value = 216.75
$ws.Range("L39").Value = 285
$ws.Range("M39").Value = 79.25
$ws.Range("N39").Value = -877
$ws.Range("H40").Value = 2073
$ws.Range("J40").Value = 2219.8
$ws.Range("L40").Value = 2219.8
$ws.Range("N40").Value = -2569.8
$ws.Range("H42").Value = 732.5
$ws.Range("I42").Value = 732.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 2197.5
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -1967.5
$ws.Range("N42").Value = ""
$ws.Range("H53").Value = 1139.9
$ws.Range("I53").Value = 1326.5294
$ws.Range("K53").Value = 1326.5294
$ws.Range("M53").Value = -689.5293999999999
$ws.Range("H54").Value = 3150
$ws.Range("I54").Value = 1916.6666
$ws.Range("J54").Value = 5000
$ws.Range("K54").Value = 1916.6666
$ws.Range("L54").Value = 5000
$ws.Range("M54").Value = -1430.6666
$ws.Range("N54").Value = -5972
$ws.Range("H55").Value = 564.2857
$ws.Range("I55").Value = 483.33334
$ws.Range("J55").Value = 625
$ws.Range("K55").Value = 483.33334
$ws.Range("L55").Value = 625
$ws.Range("M55").Value = -269.33334
$ws.Range("N55").Value = -1053
$ws.Range("H58").Value = 1645.1904
$ws.Range("I58").Value = 384.3846
$ws.Range("J58").Value = 3694
$ws.Range("K58").Value = 1153.1538
$ws.Range("L58").Value = 11082
$ws.Range("M58").Value = -1003.1538
$ws.Range("N58").Value = -11382
$ws.Range("H80").Value = 876.2222
$ws.Range("J80").Value = 660.7692
$ws.Range("L80").Value = 1982.3076
$ws.Range("N80").Value = -3978.3076
$ws.Range("H83").Value = 876.2222
$ws.Range("J83").Value = 660.7692
$ws.Range("L83").Value = 5946.922799999999
$ws.Range("N83").Value = -15930.9228
$ws.Range("H98").Value = 3014.7568
$ws.Range("I98").Value = 3245.6562
$ws.Range("J98").Value = 1537
$ws.Range("K98").Value = 3245.6562
$ws.Range("L98").Value = 1537
$ws.Range("M98").Value = -1747.6562
$ws.Range("N98").Value = -4533
$ws.Range("H122").Value = 3014.7568
$ws.Range("I122").Value = 3245.6562
$ws.Range("J122").Value = 1537
$ws.Range("K122").Value = 9736.9686
$ws.Range("L122").Value = 4611
$ws.Range("M122").Value = -7286.9686
$ws.Range("N122").Value = -9511
$ws.Range("H134").Value = 34340
$ws.Range("J134").Value = 34340
$ws.Range("L134").Value = 34340
$ws.Range("N134").Value = -44480
$ws.Range("H138").Value = 1552.2903
$ws.Range("I138").Value = 1041.7333
$ws.Range("J138").Value = 1715.234
$ws.Range("K138").Value = 3125.199900000001
$ws.Range("L138").Value = 5145.701999999999
$ws.Range("M138").Value = 2014.800099999999
$ws.Range("N138").Value = -15425.702
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4996.186
$ws.Range("I32").Value = 5270.9
$ws.Range("K32").Value = 5270.9
$ws.Range("M32").Value = -4983.9
$ws.Range("H102").Value = 18521118
$ws.Range("I102").Value = 33335614
$ws.Range("K102").Value = 33335614
$ws.Range("M102").Value = -33333992
$ws.Range("H122").Value = 2381
$ws.Range("I122").Value = 2381
$ws.Range("K122").Value = 7143
$ws.Range("M122").Value = -4693
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1871.0588
$ws.Range("I58").Value = 1524.4445
$ws.Range("J58").Value = 2261
$ws.Range("K58").Value = 1524.4445
$ws.Range("L58").Value = 2261
$ws.Range("M58").Value = -1321.4445
$ws.Range("N58").Value = -2667
$ws.Range("H95").Value = 11300.125
$ws.Range("J95").Value = 11300.125
$ws.Range("L95").Value = 11300.125
$ws.Range("N95").Value = -16792.125
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = ""
$ws.Range("H132").Value = 2046.8
$ws.Range("I132").Value = 1715.8948
$ws.Range("J132").Value = 2618.3635
$ws.Range("K132").Value = 5147.6844
$ws.Range("L132").Value = 7855.0905
$ws.Range("M132").Value = -2617.6844
$ws.Range("N132").Value = -12915.0905
$ws.Range("H136").Value = 1871.0588
$ws.Range("I136").Value = 1524.4445
$ws.Range("J136").Value = 2261
$ws.Range("K136").Value = 4573.333500000001
$ws.Range("L136").Value = 6783
$ws.Range("M136").Value = -2023.333500000001
$ws.Range("N136").Value = -11883
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15154434
$ws.Range("I131").Value = 111111430
$ws.Range("J131").Value = 3329.6316
$ws.Range("K131").Value = 333334290
$ws.Range("L131").Value = 9988.8948
$ws.Range("M131").Value = -333329250
$ws.Range("N131").Value = -20068.8948
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""
$ws.Range("H122").Value = 1884.909
$ws.Range("I122").Value = 1929.8572
$ws.Range("J122").Value = 1806.25
$ws.Range("K122").Value = 5789.571599999999
$ws.Range("L122").Value = 5418.75
$ws.Range("M122").Value = -3339.571599999999
$ws.Range("N122").Value = -10318.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1625.5555
$ws.Range("J7").Value = 2505
$ws.Range("L7").Value = 2505
$ws.Range("N7").Value = -2729
$ws.Range("H16").Value = 570.3158
$ws.Range("I16").Value = 610.4375
$ws.Range("J16").Value = 356.33334
$ws.Range("K16").Value = 610.4375
$ws.Range("L16").Value = 356.33334
$ws.Range("M16").Value = -440.4375
$ws.Range("N16").Value = -696.33334
$ws.Range("H22").Value = 641.8182
$ws.Range("I22").Value = 419
$ws.Range("K22").Value = 419
$ws.Range("M22").Value = -124
$ws.Range("H27").Value = 641.8182
$ws.Range("I27").Value = 419
$ws.Range("K27").Value = 419
$ws.Range("M27").Value = -312
$ws.Range("H93").Value = 759.44446
$ws.Range("I93").Value = 662.1429000000001
$ws.Range("K93").Value = 662.1429000000001
$ws.Range("M93").Value = 585.8570999999999
$ws.Range("H126").Value = 1625.5555
$ws.Range("J126").Value = 2505
$ws.Range("L126").Value = 7515
$ws.Range("N126").Value = -12455
$ws.Range("H136").Value = 1070.9
$ws.Range("I136").Value = 713.625
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 2140.875
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = 409.125
$ws.Range("N136").Value = -12600
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 18996.666
$ws.Range("J98").Value = 18996.666
$ws.Range("L98").Value = 18996.666
$ws.Range("N98").Value = -24986.666
$ws.Range("H100").Value = 399
$ws.Range("I100").Value = 399
$ws.Range("K100").Value = 798
$ws.Range("M100").Value = -257
$ws.Range("H109").Value = 37574.75
$ws.Range("J109").Value = 33319
$ws.Range("L109").Value = 33319
$ws.Range("N109").Value = -36093
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H126").Value = 200002800
$ws.Range("I126").Value = 250002260
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 750006780
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -750004310
$ws.Range("N126").Value = -19940
